# se modifica data para SmokeQA 29-06-2021
$wb = $excel.ActiveWorkbook

# --- DatosCuenta ---
$wsCuenta = $wb.Worksheets.Item("DatosCuenta")
$wsCuenta.Range("A2").Value = "SmokQAJuneLast"
$wsCuenta.Range("B2").Value = "SmokeNameQAJuneLast"
$wsCuenta.Range("C2").Value = 27100131
$wsCuenta.Range("D2").Value = 132

# --- DatosHogar ---
$wsHogar = $wb.Worksheets.Item("DatosHogar")
$wsHogar.Range("A2").Value = 651

# --- DatosMotor ---
$wsMotor = $wb.Worksheets.Item("DatosMotor")
$wsMotor.Range("A2").Value = "SMP033"
$wsMotor.Range("B2").Value = "ABC12SSMP033"
$wsMotor.Range("C2").Value = "ZAZ123SSMP033"

# --- DatosAP ---
$wsAP = $wb.Worksheets.Item("DatosAP")
$wsAP.Range("A2").Value = 21200132
$wsAP.Activate()
$wsAP.Range("D14").Select()

# Window geometry (maximized window) to mirror the saved view state
$win = $excel.ActiveWindow
$win.Top = -120
$win.Left = -120
$win.Width = 20730
$win.Height = 11160
